$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price): write as text (apostrophe-prefix style) so Excel
# does not auto-convert number-looking strings to actual numbers,
# then reset the cell style so no stray NumberFormat/quote-prefix style sticks.
$ws.Range("D2").Value = "'" + '70.879.35'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'" + '3.627.22'
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'" + '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'" + '605.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'" + '199.98'
$ws.Range("D6").Style = "Normal"
$ws.Range("D11").Value = "'" + '53.82'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'" + '0.0000305'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'" + '9.57'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'" + '4.203.46'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'" + '673.84'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'" + '70.960.75'
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").Value = "'" + '3.634.02'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'" + '19.02'
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Value = "'" + '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'" + '18.52'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'" + '5.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'" + '104.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Value = "'" + '10.51'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Value = "'" + '9.85'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'" + '34.14'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Value = "'" + '4.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Value = "'" + '7.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Value = "'" + '12.20'
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Value = "'" + '63.27'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = "'" + '3.935.10'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "'" + '0.0₃0864'
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Value = "'" + '516.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("D42").Value = "'" + '36.45'
$ws.Range("D42").Style = "Normal"
$ws.Range("D45").Value = "'" + '3.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'" + '3.46'
$ws.Range("D46").Style = "Normal"

# Column E (Volume 1h): plain text assignment (values are never numeric-looking)
$ws.Range("E2").Value = '  +1.84%  '
$ws.Range("E3").Value = '  +3.47%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("E6").Value = '  +1.99%  '
$ws.Range("E7").Value = '  +1.10%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("E9").Value = '  +9.19%  '
$ws.Range("E10").Value = '  +0.63%  '
$ws.Range("E11").Value = '  +1.42%  '
$ws.Range("E12").Value = '  +2.97%  '
$ws.Range("E13").Value = '  +1.52%  '
$ws.Range("E14").Value = '  +3.49%  '
$ws.Range("E15").Value = '  +12.73%  '
$ws.Range("E16").Value = '  +1.82%  '
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("E18").Value = '  +3.88%  '
$ws.Range("E19").Value = '  +0.53%  '
$ws.Range("E20").Value = '  -0.12%  '
$ws.Range("E21").Value = '  +1.90%  '
$ws.Range("E22").Value = '  +3.50%  '
$ws.Range("E23").Value = '  +2.40%  '
$ws.Range("E24").Value = '  +3.09%  '
$ws.Range("E25").Value = '  +0.11%  '
$ws.Range("E26").Value = '  -2.68%  '
$ws.Range("E27").Value = '  -2.28%  '
$ws.Range("E28").Value = '  +4.17%  '
$ws.Range("E29").Value = '  +3.63%  '
$ws.Range("E30").Value = '  +8.70%  '
$ws.Range("E31").Value = '  +3.19%  '
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("E33").Value = '  +1.63%  '
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("E35").Value = '  +5.64%  '
$ws.Range("E36").Value = '  +7.21%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("E38").Value = '  +4.64%  '
$ws.Range("E39").Value = '  -4.14%  '
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("E41").Value = '  -2.03%  '
$ws.Range("E42").Value = '  +1.39%  '
$ws.Range("E43").Value = '  +1.75%  '
$ws.Range("E44").Value = '  +2.28%  '
$ws.Range("E45").Value = '  +8.26%  '
$ws.Range("E46").Value = '  +6.32%  '
$ws.Range("E47").Value = '  +1.65%  '
$ws.Range("E48").Value = '  +2.98%  '
$ws.Range("E49").Value = '  -0.14%  '
$ws.Range("E50").Value = '  +2.15%  '
$ws.Range("E51").Value = '  +2.27%  '
